{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target edit (\"alien idea 1 save 2\"): the single paragraph of struck-through\n// text is rewritten/re-split into a new sequence of runs (same \"Idea 1 ...\"\n// theme, reworded & reordered), and the `_GoBack` bookmark that used to sit\n// right after \"Idea 1. \" moves to the very end of the paragraph (after the\n// last run, before the paragraph mark). All runs keep the original\n// <w:strike/> run formatting.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[0];\n\n// Helper: build one <w:r> with <w:strike/> formatting, preserving/omitting\n// xml:space=\"preserve\" exactly like the authored XML does (present whenever\n// the text has leading/trailing whitespace, or is empty).\nfunction run(text) {\n  const needsPreserve = text === \"\" || text !== text.trim();\n  const spaceAttr = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n  const escaped = text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n  return `<w:r><w:rPr><w:strike/></w:rPr><w:t${spaceAttr}>${escaped}</w:t></w:r>`;\n}\n\nconst newRunsXml =\n  run(\"Idea 1\") +\n  run(\"intelligent \") +\n  run(\"dinosaur\") +\n  run(\" Hollow earth. The goal is to evolve \") +\n  run(\"physical\") +\n  run(\"ly\") +\n  run(\" thru \") +\n  run(\"of the levels of the planets. \") +\n  run(\"You battle\") +\n  run(\" the humans on the surface of the planet\") +\n  run(\" using game\") +\n  run(\" \") +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  run(\"mechanies\") +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  run(\" includes guns and portals. \") +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>';\n\n// Keep the paragraph's own mark formatting/attributes intact by re-declaring\n// them verbatim on the replacement <w:p> (the paragraph range's `replace`\n// swaps the whole node, so anything not repeated here would be lost).\nconst paragraphXml =\n  '<w:p w14:paraId=\"1FB96E39\" w14:textId=\"5F40D3C0\" w:rsidR=\"007B1D5B\" w:rsidRPr=\"00EE5A4D\" w:rsidRDefault=\"00EE5A4D\">' +\n  \"<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>\" +\n  newRunsXml +\n  \"</w:p>\";\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  \"<w:body>\" +\n  paragraphXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst paragraphRange = targetParagraph.getRange();\nparagraphRange.insertOoxml(flatOpcXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument / $d is already open for us.\n#\n# Target edit (\"alien idea 1 save 2\"): the single paragraph of struck-through\n# text is rewritten/re-split into a new sequence of runs (same \"Idea 1 ...\"\n# theme, reworded & reordered), and the `_GoBack` bookmark that used to sit\n# right after \"Idea 1. \" moves to the very end of the paragraph (after the\n# last run, before the paragraph mark). All runs keep the original\n# <w:strike/> run formatting; the paragraph's own formatting (w:pPr) is left\n# untouched.\n\n$d = $word.ActiveDocument\n\n# Helper: build one <w:r> with <w:strike/> formatting, preserving/omitting\n# xml:space=\"preserve\" exactly like the authored XML does (present whenever\n# the text has leading/trailing whitespace, or is empty).\nfunction New-StrikeRun([string]$text) {\n    $needsPreserve = ($text -eq \"\") -or ($text -ne $text.Trim())\n    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n    if ($needsPreserve) {\n        return \"<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space=`\"preserve`\">$escaped</w:t></w:r>\"\n    } else {\n        return \"<w:r><w:rPr><w:strike/></w:rPr><w:t>$escaped</w:t></w:r>\"\n    }\n}\n\n$newRunsXml = (\n    (New-StrikeRun \"Idea 1\") +\n    (New-StrikeRun \"intelligent \") +\n    (New-StrikeRun \"dinosaur\") +\n    (New-StrikeRun \" Hollow earth. The goal is to evolve \") +\n    (New-StrikeRun \"physical\") +\n    (New-StrikeRun \"ly\") +\n    (New-StrikeRun \" thru \") +\n    (New-StrikeRun \"of the levels of the planets. \") +\n    (New-StrikeRun \"You battle\") +\n    (New-StrikeRun \" the humans on the surface of the planet\") +\n    (New-StrikeRun \" using game\") +\n    (New-StrikeRun \" \") +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    (New-StrikeRun \"mechanies\") +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    (New-StrikeRun \" includes guns and portals. \") +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>'\n)\n\n# Target the paragraph's range but trim off the trailing paragraph mark so\n# InsertXML only replaces the run content -- this naturally keeps the\n# paragraph's own w:pPr / w14:paraId / rsid attributes intact.\n$p = $d.Paragraphs(1)\n$r = $p.Range\n$r.MoveEnd(1, -1) | Out-Null\n\n$xml = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body><w:p>$newRunsXml</w:p></w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$r.InsertXML($xml)\n"}
